# Daily attendance processing - 2026-01-16 19:56:11
#
# In the "Recorded By" column (G), move the "System" entry to the front
# of the comma-separated recorder list on every row where it appears,
# preserving the relative order of the other recorders.
#   e.g. "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#        "backup@backdoor.com, System, system" -> "System, backup@backdoor.com, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ($text -and $text.Contains("System")) {
        $parts = $text -split ",\s*"

        $rest = @()
        $foundSystem = $false
        foreach ($p in $parts) {
            if ($p -eq "System" -and -not $foundSystem) {
                $foundSystem = $true
            } else {
                $rest += $p
            }
        }

        if ($foundSystem) {
            $newParts = @("System") + $rest
            $newText = $newParts -join ", "
            if ($newText -ne $text) {
                $cell.Value = $newText
            }
        }
    }
}
